{"js": "// Apply the tracked edits described by the diff:\n//  1. \"Hoff Meyers\" -> \"Tobias Newman Muhanguzi\" (plus the auto \"_GoBack\" bookmark\n//     that Word drops at the point of the most recent edit).\n//  2. Insert <w:proofErr> spell/grammar markers around several words that Word's\n//     proofing engine flagged: utilise, \"protection\" (as part of \"reverse current\n//     protection\"), Standardised, LoRa, Finalise, minimise.\n//\n// proofErr markers have no API surface in Word.js (they are not represented as\n// visible content), so we rebuild the OOXML of each affected paragraph and feed\n// it back with insertOoxml(..., Replace) on the paragraph's own range. That\n// keeps the paragraph's own pPr/properties (numbering, etc.) untouched because\n// we copy them through verbatim, and only the run-level text inside the\n// paragraph is restructured.\n\nasync function replaceParagraphOoxml(context, paragraph, transform) {\n  const ooxml = paragraph.getOoxml();\n  await context.sync();\n\n  const full = ooxml.value;\n  const bodyTag = \"<w:body>\";\n  const docIdx = full.indexOf(\"<w:document\");\n  const bodyOpenIdx = full.indexOf(bodyTag, docIdx) + bodyTag.length;\n  const pEndTag = \"</w:p>\";\n  const pEndIdx = full.indexOf(pEndTag, bodyOpenIdx) + pEndTag.length;\n  const paraXml = full.substring(bodyOpenIdx, pEndIdx);\n\n  const newParaXml = transform(paraXml);\n  if (newParaXml === paraXml) {\n    throw new Error(\"transform did not change paragraph xml: \" + paraXml);\n  }\n\n  const pkg =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n    \"<w:body>\" +\n    newParaXml +\n    \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\n  const range = paragraph.getRange();\n  range.insertOoxml(pkg, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nasync function findParagraphContaining(context, body, needle) {\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n  for (const p of paragraphs.items) {\n    p.load(\"text\");\n  }\n  await context.sync();\n  for (const p of paragraphs.items) {\n    if (p.text.indexOf(needle) !== -1) {\n      return p;\n    }\n  }\n  throw new Error(\"paragraph containing '\" + needle + \"' not found\");\n}\n\n// 1. \"Prepared By:\" name change, with the auto \"_GoBack\" bookmark Word leaves\n// behind at the last-edited location.\nconst nameResults = context.document.body.search(\"Hoff Meyers\", { matchCase: true });\nnameResults.load(\"items\");\nawait context.sync();\nif (nameResults.items.length === 0) {\n  throw new Error(\"'Hoff Meyers' not found\");\n}\nconst nameRange = nameResults.items[0];\nnameRange.insertText(\"Tobias Newman Muhanguzi\", Word.InsertLocation.replace);\nawait context.sync();\n\nconst newNameResults = context.document.body.search(\"Tobias Newman Muhanguzi\", { matchCase: true });\nnewNameResults.load(\"items\");\nawait context.sync();\nconst newNameEnd = newNameResults.items[0].getRange(Word.RangeLocation.end);\nnewNameEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. \"utilise\" -> wrap in spellStart/spellEnd proofErr.\nlet para = await findParagraphContaining(context, context.document.body, \"utilise\");\nawait replaceParagraphOoxml(context, para, (xml) =>\n  xml.replace(\n    \"<w:r><w:t>utilise</w:t></w:r>\",\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>utilise</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>'\n  )\n);\n\n// 3. \"reverse current protection\" -> split into two bold runs, with gramStart/\n// gramEnd proofErr wrapping the second (\"protection\").\npara = await findParagraphContaining(context, context.document.body, \"reverse current protection\");\nawait replaceParagraphOoxml(context, para, (xml) =>\n  xml.replace(\n    '<w:r w:rsidRPr=\"000F4C0D\"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>reverse current protection</w:t></w:r>',\n    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\">reverse current </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>protection</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>'\n  )\n);\n\n// 4. \"Standardised\" -> wrap in spellStart/spellEnd proofErr.\npara = await findParagraphContaining(context, context.document.body, \"Standardised\");\nawait replaceParagraphOoxml(context, para, (xml) =>\n  xml.replace(\n    \"<w:r><w:t>Standardised</w:t></w:r>\",\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Standardised</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>'\n  )\n);\n\n// 5. \"LoRa\" -> split out of the surrounding run, wrapped in spellStart/spellEnd.\npara = await findParagraphContaining(context, context.document.body, \"LoRa\");\nawait replaceParagraphOoxml(context, para, (xml) =>\n  xml.replace(\n    '<w:r w:rsidRPr=\"000F4C0D\"><w:t xml:space=\"preserve\"> Design a logic circuit to switch the FTDI UART connection between the Main MCU, GSM, and LoRa modules, allowing all to be programmed via the single USB-C port.</w:t></w:r>',\n    '<w:r><w:t xml:space=\"preserve\"> Design a logic circuit to switch the FTDI UART connection between the Main MCU, GSM, and </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      \"<w:r><w:t>LoRa</w:t></w:r>\" +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> modules, allowing all to be programmed via the single USB-C port.</w:t></w:r>'\n  )\n);\n\n// 6. \"Finalise\" and \"minimise\" -> wrap each in spellStart/spellEnd proofErr\n// (both live in the same paragraph).\n// Note: Word.js's getOoxml() silently drops the <w:lastRenderedPageBreak/>\n// hint that sits in this paragraph's first run (it's a transient layout\n// artifact, not addressable content) - put it back so the round trip via\n// getOoxml/insertOoxml doesn't lose it.\npara = await findParagraphContaining(context, context.document.body, \"Finalise\");\nawait replaceParagraphOoxml(context, para, (xml) => {\n  let out = xml.replace(\n    \"<w:rPr><w:b/><w:bCs/></w:rPr><w:t>PCB Redesign:</w:t>\",\n    \"<w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>PCB Redesign:</w:t>\"\n  );\n  out = out.replace(\n    \"<w:r><w:t>Finalise</w:t></w:r>\",\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Finalise</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>'\n  );\n  out = out.replace(\n    \"<w:r><w:t>minimise</w:t></w:r>\",\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>minimise</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>'\n  );\n  return out;\n});\n", "ps1": "# Apply the tracked edits described by the diff:\n#  1. \"Hoff Meyers\" -> \"Tobias Newman Muhanguzi\" (plus the auto \"_GoBack\" bookmark\n#     that Word drops at the point of the most recent edit).\n#  2. Insert <w:proofErr> spell/grammar markers around several words that Word's\n#     proofing engine flagged: utilise, \"protection\" (as part of \"reverse current\n#     protection\"), Standardised, LoRa, Finalise, minimise.\n#\n# proofErr markers have no surface in the Word object model (they are not\n# content, just proofing-UI markers), so we rebuild the OOXML of each affected\n# paragraph and feed it back with Range.InsertXML on the paragraph's own range.\n# That keeps the paragraph's own properties (numbering, etc.) untouched because\n# we copy them through verbatim, and only the run-level text inside the\n# paragraph is restructured.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphContaining($text) {\n  foreach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$text*\") {\n      return $p\n    }\n  }\n  throw \"paragraph containing '$text' not found\"\n}\n\nfunction Replace-ParagraphXml($paragraph, $oldStr, $newStr) {\n  $range = $paragraph.Range\n  $full = $range.WordOpenXML\n  $bodyTag = \"<w:body>\"\n  $docIdx = $full.IndexOf(\"<w:document\")\n  $bodyOpenIdx = $full.IndexOf($bodyTag, $docIdx) + $bodyTag.Length\n  $pEndTag = \"</w:p>\"\n  $pEndIdx = $full.IndexOf($pEndTag, $bodyOpenIdx) + $pEndTag.Length\n  $paraXml = $full.Substring($bodyOpenIdx, $pEndIdx - $bodyOpenIdx)\n\n  $newParaXml = $paraXml.Replace($oldStr, $newStr)\n  if ($newParaXml -eq $paraXml) {\n    throw \"replacement text not found in paragraph xml\"\n  }\n\n  $pkg = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n  $range.InsertXML($pkg)\n}\n\n# 1. \"Prepared By:\" name change, with the auto \"_GoBack\" bookmark Word leaves\n# behind at the last-edited location.\n$nameRange = $d.Content\n$nameRange.Find.Execute(\"Hoff Meyers\")\n$nameRange.Text = \"Tobias Newman Muhanguzi\"\n$nameRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $nameRange)\n\n# 2. \"utilise\" -> wrap in spellStart/spellEnd proofErr.\n$p = Get-ParagraphContaining(\"utilise\")\nReplace-ParagraphXml $p '<w:r><w:t>utilise</w:t></w:r>' '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>utilise</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>'\n\n# 3. \"reverse current protection\" -> split into two bold runs, with gramStart/\n# gramEnd proofErr wrapping the second (\"protection\").\n$p = Get-ParagraphContaining(\"reverse current protection\")\nReplace-ParagraphXml $p '<w:r w:rsidRPr=\"000F4C0D\"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>reverse current protection</w:t></w:r>' '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\">reverse current </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>protection</w:t></w:r><w:proofErr w:type=\"gramEnd\"/>'\n\n# 4. \"Standardised\" -> wrap in spellStart/spellEnd proofErr.\n$p = Get-ParagraphContaining(\"Standardised\")\nReplace-ParagraphXml $p '<w:r><w:t>Standardised</w:t></w:r>' '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Standardised</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>'\n\n# 5. \"LoRa\" -> split out of the surrounding run, wrapped in spellStart/spellEnd.\n$p = Get-ParagraphContaining(\"LoRa\")\nReplace-ParagraphXml $p '<w:r w:rsidRPr=\"000F4C0D\"><w:t xml:space=\"preserve\"> Design a logic circuit to switch the FTDI UART connection between the Main MCU, GSM, and LoRa modules, allowing all to be programmed via the single USB-C port.</w:t></w:r>' '<w:r><w:t xml:space=\"preserve\"> Design a logic circuit to switch the FTDI UART connection between the Main MCU, GSM, and </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>LoRa</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> modules, allowing all to be programmed via the single USB-C port.</w:t></w:r>'\n\n# 6. \"Finalise\" and \"minimise\" -> wrap each in spellStart/spellEnd proofErr\n# (both live in the same paragraph).\n$p = Get-ParagraphContaining(\"Finalise\")\nReplace-ParagraphXml $p '<w:r><w:t>Finalise</w:t></w:r>' '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Finalise</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>'\n$p = Get-ParagraphContaining(\"minimise\")\nReplace-ParagraphXml $p '<w:r><w:t>minimise</w:t></w:r>' '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>minimise</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>'\n"}
